$wb = $excel.ActiveWorkbook

# --- Rename the "Include from EntityNameUse" sheet to "Include #0" ---
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Name = "Include #0"

# --- Update Metadata sheet ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Bump version + regenerate date
$wsMeta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$wsMeta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" row right after "Contact" (row 10), before "Description" (row 11)
$wsMeta.Rows.Item(11).Insert()

# The freshly inserted row picked up a blank default style; re-apply the
# standard data-row formatting from the (now shifted) row below it.
$wsMeta.Range("A12:B12").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)

$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""
